{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (per commit \"updated to do list\"):\n//   - Paragraph 2 (\"//Implement auction property\") becomes a single run\n//     paragraph reading \"When paying utilities, you don't roll dice to\n//     know how much to pay\".\n//   - Paragraphs 3-6 (\"When all but one players...\", the original\n//     \"When paying utilities...\" split paragraph with the _GoBack\n//     bookmark, \"If someone is visiting jail...\", and \"//If you get to\n//     roll...\") are removed entirely.\n//   - The trailing empty paragraph keeps the _GoBack bookmark (moved\n//     there from the old \"When paying utilities...\" paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"When paying utilities, you don\\u2019t roll dice to know how much to pay\";\n\n// Find the paragraph that starts the \"//Implement auction property\" block\n// and the one that reads \"Implement full set\" (kept, unchanged) so we know\n// exactly which paragraphs in between must be removed.\nlet startIdx = -1; // \"//\" + \"Implement\" + \" auction property\"\nlet keepIdx = -1;  // \"Implement full set\"\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (startIdx === -1 && t === \"//Implement auction property\") {\n    startIdx = i;\n  }\n  if (t === \"Implement full set\") {\n    keepIdx = i;\n    break;\n  }\n}\n\nif (startIdx === -1 || keepIdx === -1 || keepIdx <= startIdx) {\n  throw new Error(\"Could not locate expected paragraphs to edit.\");\n}\n\n// Replace the text of the first paragraph in the block with the merged\n// sentence, then delete the remaining paragraphs up to (not including)\n// the \"Implement full set\" paragraph.\nparagraphs.items[startIdx].insertText(targetText, \"Replace\");\n\nfor (let i = keepIdx - 1; i > startIdx; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\n// Re-load paragraphs to get fresh references, then add the _GoBack\n// bookmark into the trailing empty paragraph (last paragraph in the body).\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs2.items[paragraphs2.items.length - 1];\nlastParagraph.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument / $d resolve to the open document.\n#\n# Target change (per commit \"updated to do list\"):\n#   - Paragraph \"//Implement auction property\" becomes a single paragraph\n#     reading \"When paying utilities, you don't roll dice to know how\n#     much to pay\".\n#   - The paragraphs \"When all but one players...\", the original\n#     \"When paying utilities...\" paragraph (which held the _GoBack\n#     bookmark), \"If someone is visiting jail...\", and \"//If you get to\n#     roll...\" are removed entirely.\n#   - The trailing empty paragraph gets the _GoBack bookmark (moved there\n#     from the old \"When paying utilities...\" paragraph).\n\n$d = $word.ActiveDocument\n\n$apos = [char]8217\n$targetText = \"When paying utilities, you don\" + $apos + \"t roll dice to know how much to pay\"\n\n# Locate the start paragraph (\"//Implement auction property\") and the\n# \"Implement full set\" paragraph that must stay immediately after the\n# edited paragraph once everything in between is removed.\n$startIndex = -1\n$keepIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)\n    if ($startIndex -eq -1 -and $t -eq \"//Implement auction property\") {\n        $startIndex = $i\n    }\n    if ($t -eq \"Implement full set\") {\n        $keepIndex = $i\n        break\n    }\n}\n\nif ($startIndex -eq -1 -or $keepIndex -eq -1 -or $keepIndex -le $startIndex) {\n    throw \"Could not locate expected paragraphs to edit.\"\n}\n\n# Delete the paragraphs strictly between the start paragraph and the kept\n# \"Implement full set\" paragraph, working backwards so earlier indices\n# stay valid.\nfor ($i = $keepIndex - 1; $i -gt $startIndex; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n# Replace the start paragraph's text (but not its paragraph mark) with the\n# merged sentence, preserving its run formatting.\n$r = $d.Paragraphs.Item($startIndex).Range\n$r.End = $r.End - 1\n$r.Text = $targetText\n\n# Add the _GoBack bookmark into the trailing empty paragraph.\n$last = $d.Paragraphs.Item($d.Paragraphs.Count)\n$d.Bookmarks.Add(\"_GoBack\", $last.Range)\n"}
